# Apply updated cryptocurrency price/volume figures to Sheet1.
#
# "Price" values that look like plain numbers (e.g. 333.24) must be forced
# to stay as text -- otherwise Excel auto-converts them to real numbers,
# which does not match the source data (some prices use "." as a thousands
# separator, e.g. 42.723.57, and the whole Price column is text). A leading
# apostrophe (`') is used to force those particular values to remain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.723.57'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '2.369.01'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = "`'333.24"
$ws.Range('E5').Value = '  +6.52%  '
$ws.Range('D6').Value = "`'101.55"
$ws.Range('E6').Value = '  -7.30%  '
$ws.Range('E7').Value = '  -0.87%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = "`'0.629"
$ws.Range('E9').Value = '  -0.83%  '
$ws.Range('D10').Value = "`'40.09"
$ws.Range('E10').Value = '  -6.73%  '
$ws.Range('E11').Value = '  -1.42%  '
$ws.Range('D12').Value = "`'8.48"
$ws.Range('E12').Value = '  -4.40%  '
$ws.Range('E13').Value = '  -3.40%  '
$ws.Range('E14').Value = '  +0.18%  '
$ws.Range('D15').Value = "`'16.49"
$ws.Range('E15').Value = '  +1.16%  '
$ws.Range('D16').Value = '2.727.93'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('D17').Value = '2.367.62'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = "`'8.11"
$ws.Range('E18').Value = '  +11.74%  '
$ws.Range('D19').Value = '42.695.24'
$ws.Range('E20').Value = '  -1.61%  '
$ws.Range('D21').Value = "`'3.79"
$ws.Range('E21').Value = '  +9.95%  '
$ws.Range('D22').Value = "`'76.72"
$ws.Range('E22').Value = '  +1.69%  '
$ws.Range('D23').Value = "`'269.05"
$ws.Range('E23').Value = '  +5.42%  '
$ws.Range('B24').Value = 'InternetComputer(DFINITY)'
$ws.Range('C24').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D24').Value = "`'10.23"
$ws.Range('E24').Value = '  +12.38%  '
$ws.Range('B25').Value = 'ImmutableX'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = "`'2.32"
$ws.Range('E25').Value = '  -10.33%  '
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('E27').Value = '  -4.17%  '
$ws.Range('D28').Value = "`'23.19"
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('E29').Value = '  -2.55%  '
$ws.Range('D30').Value = "`'176.36"
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('D31').Value = "`'3.10"
$ws.Range('E31').Value = '  -2.36%  '
$ws.Range('D32').Value = "`'0.0902"
$ws.Range('E32').Value = '  -2.59%  '
$ws.Range('D33').Value = "`'35.38"
$ws.Range('E33').Value = '  -9.74%  '
$ws.Range('D34').Value = "`'6.14"
$ws.Range('E34').Value = '  +1.36%  '
$ws.Range('E35').Value = '  -0.04%  '
$ws.Range('E36').Value = '  -6.78%  '
$ws.Range('E37').Value = '  +10.34%  '
$ws.Range('E38').Value = '  -5.03%  '
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').Value = "`'3.82"
$ws.Range('E40').Value = '  -7.61%  '
$ws.Range('E41').Value = '  +3.06%  '
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('D43').Value = "`'70.22"
$ws.Range('E43').Value = '  -3.44%  '
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').Value = "`'92.52"
$ws.Range('E45').Value = '  +32.00%  '
$ws.Range('D46').Value = "`'118.14"
$ws.Range('E46').Value = '  +6.41%  '
$ws.Range('D47').Value = "`'11.80"
$ws.Range('E47').Value = '  -7.64%  '
$ws.Range('E48').Value = '  -2.10%  '
$ws.Range('D49').Value = "`'9.21"
$ws.Range('E49').Value = '  -1.20%  '
$ws.Range('E50').Value = '  -2.77%  '
$ws.Range('D51').Value = '1.568.16'
$ws.Range('E51').Value = '  +4.96%  '
